$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2460
$ws.Range("E2").Value = 259
$ws.Range("F2").Value = 259
$ws.Range("G2").Value = 191
$ws.Range("H2").Value = 130
$ws.Range("I2").Value = 122
$ws.Range("J2").Value = 7
$ws.Range("K2").Value = 4217
$ws.Range("L2").Value = 2625
$ws.Range("M2").Value = 1592
$ws.Range("N2").Value = 1580
$ws.Range("O2").Value = 12
$ws.Range("P2").Value = 468
$ws.Range("Q2").Value = 293
$ws.Range("R2").Value = -376
$ws.Range("S2").Value = -26
$ws.Range("T2").Value = 379
$ws.Range("U2").Value = -87
$ws.Range("V2").Value = 1740
$ws.Range("W2").Value = 10.53
$ws.Range("X2").Value = 5.27
$ws.Range("Y2").Value = 8.31
$ws.Range("Z2").Value = 3.22
$ws.Range("AA2").Value = 164.89
$ws.Range("AB2").Value = 239.88
$ws.Range("AC2").Value = 1324
$ws.Range("AD2").Value = 10.99
$ws.Range("AE2").Value = 16905
$ws.Range("AF2").Value = 0.86
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 1.03
$ws.Range("AI2").Value = 11.48
$ws.Range("AJ2").Value = 9352731

# Row 3
$ws.Range("D3").Value = 2677
$ws.Range("E3").Value = 185
$ws.Range("F3").Value = 185
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 25
$ws.Range("I3").Value = 17
$ws.Range("J3").Value = 8
$ws.Range("K3").Value = 4676
$ws.Range("L3").Value = 3072
$ws.Range("M3").Value = 1604
$ws.Range("N3").Value = 1585
$ws.Range("O3").Value = 20
$ws.Range("P3").Value = 469
$ws.Range("Q3").Value = 275
$ws.Range("R3").Value = -397
$ws.Range("S3").Value = 66
$ws.Range("T3").Value = 390
$ws.Range("U3").Value = -115
$ws.Range("V3").Value = 1993
$ws.Range("W3").Value = 6.92
$ws.Range("X3").Value = 0.94
$ws.Range("Y3").Value = 1.07
$ws.Range("Z3").Value = 0.56
$ws.Range("AA3").Value = 191.46
$ws.Range("AB3").Value = 238.43
$ws.Range("AC3").Value = 180
$ws.Range("AD3").Value = 70.17
$ws.Range("AE3").Value = 16914
$ws.Range("AF3").Value = 0.75
$ws.Range("AG3").Value = 150
$ws.Range("AH3").Value = 1.19
$ws.Range("AI3").Value = 83.33
$ws.Range("AJ3").Value = 9373525

# Row 4
$ws.Range("D4").Value = 2903
$ws.Range("E4").Value = 178
$ws.Range("F4").Value = 178
$ws.Range("G4").Value = 27
$ws.Range("H4").Value = 5
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = -2
$ws.Range("K4").Value = 4714
$ws.Range("L4").Value = 3093
$ws.Range("M4").Value = 1621
$ws.Range("N4").Value = 1603
$ws.Range("O4").Value = 18
$ws.Range("P4").Value = 469
$ws.Range("Q4").Value = 356
$ws.Range("R4").Value = -208
$ws.Range("S4").Value = -100
$ws.Range("T4").Value = 210
$ws.Range("U4").Value = 146
$ws.Range("V4").Value = 1998
$ws.Range("W4").Value = 6.13
$ws.Range("X4").Value = 0.17
$ws.Range("Y4").Value = 0.4
$ws.Range("Z4").Value = 0.1
$ws.Range("AA4").Value = 190.77
$ws.Range("AB4").Value = 236.01
$ws.Range("AC4").Value = 68
$ws.Range("AD4").Value = 166.2
$ws.Range("AE4").Value = 17113
$ws.Range("AF4").Value = 0.66
$ws.Range("AG4").Value = 150
$ws.Range("AH4").Value = 1.32
$ws.Range("AI4").Value = 219.53
$ws.Range("AJ4").Value = 9373525

# Row 5
$ws.Range("D5").Value = 2987
$ws.Range("E5").Value = 160
$ws.Range("F5").Value = 160
$ws.Range("G5").Value = 99
$ws.Range("H5").Value = 54
$ws.Range("I5").Value = 53
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 4393
$ws.Range("L5").Value = 2859
$ws.Range("M5").Value = 1533
$ws.Range("N5").Value = 1523
$ws.Range("O5").Value = 10
$ws.Range("P5").Value = 469
$ws.Range("Q5").Value = 416
$ws.Range("R5").Value = -52
$ws.Range("S5").Value = -438
$ws.Range("T5").Value = 98
$ws.Range("U5").Value = 318
$ws.Range("V5").Value = 1599
$ws.Range("W5").Value = 5.36
$ws.Range("X5").Value = 1.82
$ws.Range("Y5").Value = 3.42
$ws.Range("Z5").Value = 1.2
$ws.Range("AA5").Value = 186.47
$ws.Range("AB5").Value = 244.97
$ws.Range("AC5").Value = 570
$ws.Range("AD5").Value = 16.72
$ws.Range("AE5").Value = 16514
$ws.Range("AF5").Value = 0.58
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 2.63
$ws.Range("AI5").Value = 43.33
$ws.Range("AJ5").Value = 9373525

# Row 6
$ws.Range("D6").Value = 3396
$ws.Range("E6").Value = 237
$ws.Range("F6").Value = 237
$ws.Range("G6").Value = 120
$ws.Range("H6").Value = 65
$ws.Range("I6").Value = 62
$ws.Range("K6").Value = 4342
$ws.Range("L6").Value = 2537
$ws.Range("M6").Value = 1805
$ws.Range("N6").Value = 1791
$ws.Range("P6").Value = 548
$ws.Range("Q6").Value = 337
$ws.Range("R6").Value = -81
$ws.Range("S6").Value = -182
$ws.Range("T6").Value = 67
$ws.Range("U6").Value = 270
$ws.Range("V6").Value = 1410
$ws.Range("W6").Value = 6.98
$ws.Range("X6").Value = 1.92
$ws.Range("Y6").Value = 3.73
$ws.Range("Z6").Value = 1.49
$ws.Range("AA6").Value = 140.57
$ws.Range("AB6").Value = 237.49
$ws.Range("AC6").Value = 637
$ws.Range("AD6").Value = 16.1
$ws.Range("AE6").Value = 16661
$ws.Range("AF6").Value = 0.62
$ws.Range("AG6").Value = 200
$ws.Range("AH6").Value = 1.95
$ws.Range("AI6").Value = 34.79
$ws.Range("AJ6").Value = 10952635

# Clear rows 7-9 (D:AJ) - data no longer reported for these periods
$ws.Range("D7:AJ9").ClearContents()
